$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "balance" column header
$ws.Range("E1").Value = "balance"

# Add balance values for each row
$ws.Range("E2").Value = 1234.34
$ws.Range("E3").Value = 1212.23
$ws.Range("E4").Value = 4343.02
$ws.Range("E5").Value = 344.94

# Move selection to reflect the final active cell after data entry
$ws.Range("E7").Select()
